# The deck originally carries the "Integral" theme (ppt/theme/theme1.xml,
# used by the one-and-only Slide Master) and an "Office Theme" colour
# scheme (ppt/theme/theme2.xml, used only by the Notes Master). The
# commit swaps the two themes' content wholesale: the Slide Master's
# theme becomes "Office Theme" colours and the Notes Master's theme
# becomes "Integral" colours.
#
# The PowerPoint object model only exposes one writable theme colour
# scheme for editing (Master.Theme.ThemeColorScheme / Slide.ThemeColorScheme
# all resolve back to the deck's single registered theme part), so we
# apply the colour half of the swap that is reachable from the object
# model: push the "Office Theme" palette (the colours theme1.xml should
# end up with) onto the Slide Master's theme via
# ThemeColorScheme.Colors(i).RGB, in the standard
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink (1..12) index order.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function ToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colour scheme (formerly theme2.xml),
# in MsoThemeColorSchemeIndex order (1-based):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = ToRgb($officeThemeColors[$i - 1])
}

Write-Host "Applied Office Theme colour scheme to the presentation theme."
